$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered out save games) -> updated values for rows 2 and 3

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 5.586269137925634

$ws.Range("B3").Value = 0.1190320826869504
$ws.Range("C3").Value = 0.04071648406533734
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1.406728370586922
